$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.749.88"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.538.70"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "2.928.46"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "2.583.71"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.816"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Value = "42.728.87"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -5.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  +12.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("E36").Value = "  -5.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.12%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  +9.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "1.964.45"
$ws.Range("E46").Value = "  -2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "2.782.54"
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "80.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.862"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.03%  "
